# Applies the "Updated cryptos list" data refresh to sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) cells whose new text happens to look numeric ---
# Force the whole Price column to Text format first so assigning these
# strings does not get silently coerced into floating point numbers
# (which would lose trailing zeros / exact formatting, e.g. "7.500").
$priceRange = $ws.Range("D5:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D5").Value = "240.06"
$ws.Range("D6").Value = "0.6284"
$ws.Range("D7").Value = "0.9996"
$ws.Range("D8").Value = "0.07633"
$ws.Range("D9").Value = "0.2905"
$ws.Range("D10").Value = "24.73"
$ws.Range("D11").Value = "0.07738"
$ws.Range("D13").Value = "0.6789"
$ws.Range("D15").Value = "83.24"
$ws.Range("D16").Value = "6.152"
$ws.Range("D18").Value = "227.53"
$ws.Range("D21").Value = "7.500"
$ws.Range("D23").Value = "158.55"
$ws.Range("D24").Value = "0.1383"
$ws.Range("D25").Value = "8.399"
$ws.Range("D27").Value = "1.374"
$ws.Range("D29").Value = "0.05598"
$ws.Range("D30").Value = "4.114"
$ws.Range("D31").Value = "4.075"
$ws.Range("D34").Value = "0.6942"
$ws.Range("D35").Value = "2.579"
$ws.Range("D39").Value = "6.389"
$ws.Range("D40").Value = "0.9035"
$ws.Range("D42").Value = "101.47"
$ws.Range("D43").Value = "66.05"
$ws.Range("D44").Value = "7.178"
$ws.Range("D45").Value = "0.4009"
$ws.Range("D46").Value = "8.999"
$ws.Range("D47").Value = "1.672"
$ws.Range("D48").Value = "0.1137"
$ws.Range("D49").Value = "0.05702"
$ws.Range("D50").Value = "0.4628"
$ws.Range("D51").Value = "2.531"

# Restore the default (un-styled) cell style now that the text values are
# committed, so the cells keep matching their original (style-less) look.
$priceRange.Style = "Normal"

# --- Remaining cells: plain text, not number-like, no special handling needed ---
$ws.Range("D2").Value = "29.375.45"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.847.79"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "29.395.57"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("E27").Value = "  +4.95%  "
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("E34").Value = "  -2.17%  "
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").Value = "1.230.18"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("E38").Value = "  -2.29%  "
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("E51").Value = "  -0.17%  "

